$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update cell text/values to their new contents per the target layout.
$ws.Range("B10").Value = "3268262 - Carlos Renato Menegatti"
$ws.Range("C10").Value = "3268262 - Carlos Renato Menegatti"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Field Electrostatic Equipotential Mapping; Introduction to Direct Current Circuits; Resistance, Resistivity and Electrical Current; Kirchoff Laws; Capacitors; Voltmeters, Ammeters and Ohmmeters; Oscilloscopes; Magnetostatic Field; Faraday's Induction Law; Circuits RL and RC;"
$ws.Range("C14").Value = "Field Electrostatic Equipotential Mapping; Introduction to Direct Current Circuits; Resistance, Resistivity and Electrical Current; Kirchoff Laws; Capacitors; Voltmeters, Ammeters and Ohmmeters; Oscilloscopes; Magnetostatic Field; Faraday's Induction Law; Circuits RL and RC;"
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1) Electrostatic Field and Equipotential Mapping: Parallel plates Field, A point charge Field, insulating effect and conductor.2) Ohm’s Law: ohmic resistors, resistors non-ohmic.3) Resistance and Electric current: Ohm's Law, Drude model.4) Direct Current Circuits: Kirchoff laws.5) Capacitors: Capacitors association, load and discharge a capacitor.6) Voltmeters, Ammeters and ohmmeters: Galvanometer operation principle, Voltmeters Construction, Ammeters and ohmmeters.7) Oscilloscope: Oscilloscope Operation Principle.8) Magnetostatic Field: Biot-Savart law, Ampere's law, Hall effect.9) Faraday's Law of Induction: Mutual inductance and self-inductance, AC voltage generation.10) RL and RC in DC circuits"
$ws.Range("C16").Value = "1) Electrostatic Field and Equipotential Mapping: Parallel plates Field, A point charge Field, insulating effect and conductor.2) Ohm’s Law: ohmic resistors, resistors non-ohmic.3) Resistance and Electric current: Ohm's Law, Drude model.4) Direct Current Circuits: Kirchoff laws.5) Capacitors: Capacitors association, load and discharge a capacitor.6) Voltmeters, Ammeters and ohmmeters: Galvanometer operation principle, Voltmeters Construction, Ammeters and ohmmeters.7) Oscilloscope: Oscilloscope Operation Principle.8) Magnetostatic Field: Biot-Savart law, Ampere's law, Hall effect.9) Faraday's Law of Induction: Mutual inductance and self-inductance, AC voltage generation.10) RL and RC in DC circuits"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "3268262 - Carlos Renato Menegatti"
$ws.Range("C18").Value = "3268262 - Carlos Renato Menegatti"
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B23").Value = "LOB1038 -  Física Experimental I  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOB1038 -  Física Experimental I  (Requisito fraco)`n"

# 2) Clear cell A23 leftover (old "Requisitos:" row becomes row 22; row 23 only has B/C).
$ws.Range("A23").ClearContents()

# 3) Remove the now-duplicated trailing row (old row 24), shrinking the sheet to A1:C23.
$ws.Rows.Item(24).Delete()

# 4) Fix up row heights so they match the new layout exactly.
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30
